# 2ИСИП-422_ДисМат_.xlsx — AutoCommit_21 декабря 2023 г. 14:15:36_SibNout2023
#
# Adds "авансом" (S/T columns = grade given in advance) and "отчислена"
# (expelled) notes for a handful of students, and moves the sheet's
# scroll/selection state so S32 (Чекаева Софья, marked "отчислена") is
# the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Students who received credit "in advance" (авансом): S = points, T = note.
$ws.Range("S5").Value  = 5
$ws.Range("T5").Value  = "авансом"

$ws.Range("S10").Value = 5
$ws.Range("T10").Value = "авансом"

$ws.Range("S19").Value = 5
$ws.Range("T19").Value = "авансом"

$ws.Range("S25").Value = 5
$ws.Range("T25").Value = "авансом"

$ws.Range("S27").Value = 3
$ws.Range("T27").Value = "авансом"

$ws.Range("S29").Value = 5
$ws.Range("T29").Value = "авансом"

# Чекаева Софья — expelled.
$ws.Range("S31").Value = "отчислена"

# Move the view so the newly-edited cell is in focus.
$ws.Range("S32").Select()
